# Changed fuel types to string ID.
# Column A ("id") previously held a plain numeric row index (1..11). It is
# changed to hold a short, stable string identifier for each fuel type so
# downstream code can reference fuels by name instead of position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row -> new string id, in sheet order (row 2 = Electricity ... row 12 = District Heat, Hot Water)
$fuelIds = @{
    2  = "elec"
    3  = "ng"
    4  = "propane"
    5  = "oil1"
    6  = "oil2"
    7  = "birch"
    8  = "spruce"
    9  = "pellets"
    10 = "coal"
    11 = "steam"
    12 = "hot_water"
}

foreach ($row in $fuelIds.Keys) {
    $ws.Range("A$row").Value = $fuelIds[$row]
}
